$d = $word.ActiveDocument

$replacements = @(
    @("2024-11-26 Tuesday", "2024-11-27 Wednesday"),
    @("923÷8=", "761÷4="),
    @("925÷2=", "976÷9="),
    @("797÷4=", "939÷5="),
    @("407÷5=", "197÷3="),
    @("422÷7=", "409÷4="),
    @("826÷7=", "740÷5="),
    @("813÷5=", "827÷8="),
    @("159÷8=", "183÷9="),
    @("700÷6=", "308÷5="),
    @("809÷9=", "875÷4="),
    @("825÷6=", "708÷6="),
    @("245÷9=", "394÷5="),
    @("686÷6=", "280÷4="),
    @("109÷2=", "306÷9="),
    @("642÷2=", "959÷3="),
    @("114÷4=", "833÷3="),
    @("499÷9=", "556÷8="),
    @("566÷3=", "315÷2="),
    @("899÷9=", "446÷5="),
    @("937÷6=", "214÷7="),
    @("670÷7=", "885÷6="),
    @("275÷2=", "926÷6="),
    @("808÷9=", "129÷2="),
    @("803÷6=", "249÷7="),
    @("726÷6=", "355÷9=")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
